# Re-run SGNN to annotate dialog acts following clean up work to the
# original transcripts. Updates the DAMSLTag (col I) and DialogAct (col J)
# columns for the rows whose automatic dialog-act classification changed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = "aa"
$ws.Range("J2").Value = "Agree/Accept"

$ws.Range("I5").Value = "aa"
$ws.Range("J5").Value = "Agree/Accept"

$ws.Range("I14").Value = "sv"
$ws.Range("J14").Value = "Statement-opinion"

$ws.Range("I15").Value = "sd"
$ws.Range("J15").Value = "Statement-non-opinion"

$ws.Range("I36").Value = "sv"
$ws.Range("J36").Value = "Statement-opinion"

$ws.Range("I42").Value = "sd"
$ws.Range("J42").Value = "Statement-non-opinion"

$ws.Range("I48").Value = "b"
$ws.Range("J48").Value = "Acknowledge (Backchannel)"

$ws.Range("I51").Value = "sv"
$ws.Range("J51").Value = "Statement-opinion"

$ws.Range("I53").Value = "sv"
$ws.Range("J53").Value = "Statement-opinion"

$ws.Range("I58").Value = "aa"
$ws.Range("J58").Value = "Agree/Accept"

$ws.Range("I64").Value = "sd"
$ws.Range("J64").Value = "Statement-non-opinion"

$ws.Range("I72").Value = "sd"
$ws.Range("J72").Value = "Statement-non-opinion"

$ws.Range("I78").Value = "aa"
$ws.Range("J78").Value = "Agree/Accept"

$ws.Range("I91").Value = "sv"
$ws.Range("J91").Value = "Statement-opinion"

$ws.Range("I92").Value = "sd"
$ws.Range("J92").Value = "Statement-non-opinion"

$ws.Range("I98").Value = "%"
$ws.Range("J98").Value = "Uninterpretable"

$ws.Range("I103").Value = "%"
$ws.Range("J103").Value = "Uninterpretable"

$ws.Range("I105").Value = "sv"
$ws.Range("J105").Value = "Statement-opinion"
